$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("English")

# --- Rows 109-150: Key/Value pairs (row 109-110 rename existing keys, 111-150 are new) ---
$ws.Range("A109").Value = 'ReportUserOption1Txt'
$ws.Range("B109").Value = '1. Bullying, Harassing or Intimidating'
$ws.Range("A110").Value = 'ReportUserOption2Txt'
$ws.Range("B110").Value = '2. Shared my or someone else''s private information.'
$ws.Range("A111").Value = 'ReportUserOption3Txt'
$ws.Range("B111").Value = '3. Spam, Scam or Phishing'
$ws.Range("A112").Value = 'ReportUserOption4Txt'
$ws.Range("B112").Value = '4. Intellectual Property Infringement'
$ws.Range("A113").Value = 'ReportUserOption5Txt'
$ws.Range("B113").Value = '5. Self Injury or Suicide'
$ws.Range("A114").Value = 'ReportUserOption6Txt'
$ws.Range("B114").Value = '6. Conducts Illegal Activities'
$ws.Range("A115").Value = 'ReportUserOption7Txt'
$ws.Range("B115").Value = '7. Pretending to be me or someone else'
$ws.Range("A116").Value = 'ReportUserOption8Txt'
$ws.Range("B116").Value = '8. Something Else'
$ws.Range("A117").Value = 'ReportUserOption9Txt'
$ws.Range("B117").Value = '9. Not Interested/ Not Relevant'
$ws.Range("A118").Value = 'ReportKooOption1Txt'
$ws.Range("B118").Value = '1. Hate Speech and Discrimination'
$ws.Range("A119").Value = 'ReportKooOption2Txt'
$ws.Range("B119").Value = '2. Religiously Offensive Content'
$ws.Range("A120").Value = 'ReportKooOption3Txt'
$ws.Range("B120").Value = '3. Terrorism and Extremism'
$ws.Range("A121").Value = 'ReportKooOption4Txt'
$ws.Range("B121").Value = '4. Dangerous, Violent Content (Death or Injury)'
$ws.Range("A122").Value = 'ReportKooOption5Txt'
$ws.Range("B122").Value = '5. Graphic, Obscene or Sexual Content'
$ws.Range("A123").Value = 'ReportKooOption6Txt'
$ws.Range("B123").Value = '6. Provoking Koo / Comment'
$ws.Range("A124").Value = 'ReportKooOption7Txt'
$ws.Range("B124").Value = '7. Misinformation and Disinformation'
$ws.Range("A125").Value = 'ReportKooOption8Txt'
$ws.Range("B125").Value = '8. Malicious Programmes, URL or Code'
$ws.Range("A126").Value = 'ReportKooOption9Txt'
$ws.Range("B126").Value = '9. Child Abuse'
$ws.Range("A127").Value = 'ReportKooOption10Txt'
$ws.Range("B127").Value = '10. Not Interested/ Not Relevant'
$ws.Range("A128").Value = 'OnBoardingPage1HeaderText'
$ws.Range("B128").Value = 'Add your name'
$ws.Range("A129").Value = 'OnBoardingPage2HeaderText'
$ws.Range("B129").Value = 'Add Profile Picture'
$ws.Range("A130").Value = 'OnBoardingPage3HeaderText'
$ws.Range("B130").Value = 'Your Preferences'
$ws.Range("A131").Value = 'OnBoardingPage4HeaderText'
$ws.Range("B131").Value = 'People you can follow'
$ws.Range("A132").Value = 'OnBoardingPage5HeaderText'
$ws.Range("B132").Value = 'Follow Topics'
$ws.Range("A133").Value = 'AccountSectionOption1Text'
$ws.Range("B133").Value = 'Edit Profile'
$ws.Range("A134").Value = 'AccountSectionOption2Text'
$ws.Range("B134").Value = 'Account Information'
$ws.Range("A135").Value = 'AccountSectionOption3Text'
$ws.Range("B135").Value = 'Language'
$ws.Range("A136").Value = 'AccountSectionOption4Text'
$ws.Range("B136").Value = 'Theme'
$ws.Range("A137").Value = 'AccountSectionOption5Text'
$ws.Range("B137").Value = 'Share your Koo profile'
$ws.Range("A138").Value = 'AccountSectionOption6Text'
$ws.Range("B138").Value = 'Manage Subscriptions'
$ws.Range("A139").Value = 'AccountSectionOption7Text'
$ws.Range("B139").Value = 'Manage Blocked Users'
$ws.Range("A140").Value = 'AccountSectionOption8Text'
$ws.Range("B140").Value = 'Apply for Self Verification'
$ws.Range("A141").Value = 'AccountSectionOption9Text'
$ws.Range("B141").Value = 'Apply for Eminence'
$ws.Range("A142").Value = 'AccountSectionOption10Text'
$ws.Range("B142").Value = 'Migrate From Twitter'
$ws.Range("A143").Value = 'AccountSectionOption11Text'
$ws.Range("B143").Value = 'Logout'
$ws.Range("A144").Value = 'AccountSectionOption12Text'
$ws.Range("B144").Value = 'Delete'
$ws.Range("A145").Value = 'AccountSectionOption1TextById'
$ws.Range("B145").Value = 'Add Account'
$ws.Range("A146").Value = 'AccountSectionOption2TextById'
$ws.Range("B146").Value = 'Manage Account'
$ws.Range("A147").Value = 'SearchAllResultTabs1'
$ws.Range("B147").Value = 'All'
$ws.Range("A148").Value = 'SearchAllResultTabs2'
$ws.Range("B148").Value = 'People'
$ws.Range("A149").Value = 'SearchAllResultTabs3'
$ws.Range("B149").Value = 'Hashtags'
$ws.Range("A150").Value = 'SearchAllResultTabs4'
$ws.Range("B150").Value = 'Koos'

# --- Formatting for the newly written rows, matching the sheet-wide convention ---
# Column A: left aligned, no wrap. Column B: left aligned, wrap text.
$colA = $ws.Range("A109:A150")
$colA.HorizontalAlignment = -4131
$colB = $ws.Range("B109:B150")
$colB.HorizontalAlignment = -4131
$colB.WrapText = $true

# --- Rows 133-146 (Account section options) were pasted in from another sheet and carry
# --- their own explicit font/theme-color stamp distinct from the rest of the sheet. ---
$acctA = $ws.Range("A133:A146")
$acctA.Font.ThemeColor = 1
$acctA.Font.TintAndShade = 0
$acctA.HorizontalAlignment = -4131

$acctB = $ws.Range("B133:B146")
$acctB.Font.ThemeColor = 1
$acctB.Font.TintAndShade = 0
$acctB.HorizontalAlignment = -4131
$acctB.WrapText = $true

# --- Make "English" the active sheet/tab and leave the selection just past the new data ---
$ws.Activate()
$ws.Range("A151").Select()
